$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark finished topics as done by clearing their cells: all of Psych
# plus half of Gen Chem/Orgo (Equilibrium, Titration Curves, Spectroscopy,
# Lab Techniques).
$cellsToClear = @("B1", "B2", "C2", "E2", "C3", "C4", "E4", "C5", "C6", "C7", "C8")
foreach ($cellRef in $cellsToClear) {
    $ws.Range($cellRef).ClearContents()
}

# Move the active selection to E12
$ws.Range("E12").Select()
